# Doing Updates for Financials
# Insert a new first data column (D) on the BANF sheet, shifting the
# existing year columns D:K one column to the right (now E:L), then
# populate the new column D with the latest year's (2018-12-31) figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BANF")

# Insert a new column before column D; this shifts D:K -> E:L.
$ws.Range("D1").EntireColumn.Insert()

# The inserted column does not automatically pick up the number
# formatting of its neighbours, so copy it over from column E (the
# former column D) for each of the three financial-statement blocks.
$ws.Range("E7:E35").Copy()
$ws.Range("D7:D35").PasteSpecial(-4122)
$ws.Range("E38:E77").Copy()
$ws.Range("D38:D77").PasteSpecial(-4122)
$ws.Range("E80:E102").Copy()
$ws.Range("D80:D102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Income Statement (new 2018-12-31 column) ---
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 303200
$ws.Range("D9").Value = "NA"
$ws.Range("D10").Value = "NA"
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("D15").Value = -13500
$ws.Range("D17").Value = 46500
$ws.Range("D18").Value = 256700
$ws.Range("D20").Value = -96900
$ws.Range("D21").Value = 173300
$ws.Range("D22").Value = 0
$ws.Range("D23").Value = 159800
$ws.Range("D24").Value = 33900
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 125800
$ws.Range("D27").Value = 125800
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = 96900
$ws.Range("D33").Value = 125800
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 125800

# --- Balance Sheet (new 2018-12-31 column) ---
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 228400
$ws.Range("D42").Value = 1195800
$ws.Range("D43").Value = 0
$ws.Range("D44").Value = 0
$ws.Range("D45").Value = 0
$ws.Range("D46").Value = 0
$ws.Range("D47").Value = 0
$ws.Range("D48").Value = 174400
$ws.Range("D49").Value = 96200
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 0
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 7574300
$ws.Range("D57").Value = 37500
$ws.Range("D58").Value = 0
$ws.Range("D59").Value = 0
$ws.Range("D60").Value = 0
$ws.Range("D61").Value = 26800
$ws.Range("D62").Value = 0
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 6671500
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 722600
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 902800
$ws.Range("D77").Value = 0

# --- Cash Flow Statement (new 2018-12-31 column) ---
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = 125800
$ws.Range("D83").Value = 13500
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 131100
$ws.Range("D91").Value = -51900
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -284900
$ws.Range("D96").Value = -30300
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = -179800
$ws.Range("D101").Value = 0
$ws.Range("D102").Value = -333600
